# Scheduled runner: refresh Leve vendor/market profit figures (H-N) across job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 100
$ws.Cells.Item(100, 8).Value = 1846
$ws.Cells.Item(100, 9).Value = 1975.2
$ws.Cells.Item(100, 10).Value = 1200
$ws.Cells.Item(100, 11).Value = 1975.2
$ws.Cells.Item(100, 12).Value = 1200
$ws.Cells.Item(100, 13).Value = -1434.2
$ws.Cells.Item(100, 14).Value = -2282

# ALC row 137
$ws.Cells.Item(137, 8).Value = 1762130.9
$ws.Cells.Item(137, 9).Value = 5730.927
$ws.Cells.Item(137, 10).Value = 6262906
$ws.Cells.Item(137, 11).Value = 17192.781
$ws.Cells.Item(137, 12).Value = 18788718
$ws.Cells.Item(137, 13).Value = -14642.781
$ws.Cells.Item(137, 14).Value = -18793818

$ws = $wb.Worksheets.Item("ARM")
# ARM row 43
$ws.Cells.Item(43, 8).Value = 44040.2
$ws.Cells.Item(43, 9).Value = 44748.5
$ws.Cells.Item(43, 10).Value = 43568
$ws.Cells.Item(43, 11).Value = 44748.5
$ws.Cells.Item(43, 12).Value = 43568
$ws.Cells.Item(43, 13).Value = -44435.5
$ws.Cells.Item(43, 14).Value = -44194

# ARM row 46
$ws.Cells.Item(46, 8).Value = 5150.5
$ws.Cells.Item(46, 10).Value = 6200.6665
$ws.Cells.Item(46, 12).Value = 6200.6665
$ws.Cells.Item(46, 14).Value = -6838.6665

# ARM row 52
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 13).ClearContents()

# ARM row 61
$ws.Cells.Item(61, 8).Value = 956291.1
$ws.Cells.Item(61, 9).Value = 28675.477
$ws.Cells.Item(61, 11).Value = 28675.477
$ws.Cells.Item(61, 13).Value = -28463.477

# ARM row 88
$ws.Cells.Item(88, 8).Value = 2966.3333
$ws.Cells.Item(88, 10).Value = 2966.3333
$ws.Cells.Item(88, 12).Value = 2966.3333
$ws.Cells.Item(88, 14).Value = -3778.3333

# ARM row 91
$ws.Cells.Item(91, 8).Value = 2966.3333
$ws.Cells.Item(91, 10).Value = 2966.3333
$ws.Cells.Item(91, 12).Value = 2966.3333
$ws.Cells.Item(91, 14).Value = -5774.3333

# ARM row 110
$ws.Cells.Item(110, 8).Value = 1624.8823
$ws.Cells.Item(110, 9).Value = 1615
$ws.Cells.Item(110, 11).Value = 1615
$ws.Cells.Item(110, 13).Value = 430

# ARM row 122
$ws.Cells.Item(122, 8).Value = 2599.6667
$ws.Cells.Item(122, 9).Value = 2599.6667
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 7799.000100000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -5349.000100000001
$ws.Cells.Item(122, 14).ClearContents()

# ARM row 136
$ws.Cells.Item(136, 8).Value = 956291.1
$ws.Cells.Item(136, 9).Value = 28675.477
$ws.Cells.Item(136, 11).Value = 86026.431
$ws.Cells.Item(136, 13).Value = -83476.431

$ws = $wb.Worksheets.Item("BSM")
# BSM row 14
$ws.Cells.Item(14, 8).Value = 11000
$ws.Cells.Item(14, 9).Value = 11000
$ws.Cells.Item(14, 11).Value = 11000
$ws.Cells.Item(14, 13).Value = -10828

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Cells.Item(22, 8).Value = 575.3226
$ws.Cells.Item(22, 9).Value = 560.8077
$ws.Cells.Item(22, 10).Value = 650.8
$ws.Cells.Item(22, 11).Value = 560.8077
$ws.Cells.Item(22, 12).Value = 650.8
$ws.Cells.Item(22, 13).Value = -210.8077
$ws.Cells.Item(22, 14).Value = -1350.8

# CRP row 86
$ws.Cells.Item(86, 8).Value = 20088.666
$ws.Cells.Item(86, 9).Value = 8268
$ws.Cells.Item(86, 11).Value = 8268
$ws.Cells.Item(86, 13).Value = -7145

# CRP row 89
$ws.Cells.Item(89, 8).Value = 20088.666
$ws.Cells.Item(89, 9).Value = 8268
$ws.Cells.Item(89, 11).Value = 41340
$ws.Cells.Item(89, 13).Value = -35724

# CRP row 122
$ws.Cells.Item(122, 8).Value = 2241.8333
$ws.Cells.Item(122, 9).Value = 2241.8333
$ws.Cells.Item(122, 11).Value = 6725.499899999999
$ws.Cells.Item(122, 13).Value = -4275.499899999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Cells.Item(5, 8).Value = 1377.8
$ws.Cells.Item(5, 9).Value = 1086.1
$ws.Cells.Item(5, 11).Value = 3258.3
$ws.Cells.Item(5, 13).Value = -3146.3

# CUL row 87
$ws.Cells.Item(87, 8).Value = 9504.666999999999
$ws.Cells.Item(87, 9).Value = 1390.1666
$ws.Cells.Item(87, 11).Value = 4170.4998
$ws.Cells.Item(87, 13).Value = -2922.4998

# CUL row 90
$ws.Cells.Item(90, 8).Value = 9504.666999999999
$ws.Cells.Item(90, 9).Value = 1390.1666
$ws.Cells.Item(90, 11).Value = 12511.4994
$ws.Cells.Item(90, 13).Value = -6271.499400000001

# CUL row 106
$ws.Cells.Item(106, 8).Value = 10151
$ws.Cells.Item(106, 10).Value = 10151
$ws.Cells.Item(106, 12).Value = 30453
$ws.Cells.Item(106, 14).Value = -32345

# CUL row 135
$ws.Cells.Item(135, 8).Value = 1377.8
$ws.Cells.Item(135, 9).Value = 1086.1
$ws.Cells.Item(135, 11).Value = 9774.9
$ws.Cells.Item(135, 13).Value = -7239.9

# CUL row 139
$ws.Cells.Item(139, 8).Value = 9617974
$ws.Cells.Item(139, 9).Value = 15626723
$ws.Cells.Item(139, 10).Value = 3975.2
$ws.Cells.Item(139, 11).Value = 46880169
$ws.Cells.Item(139, 12).Value = 11925.6
$ws.Cells.Item(139, 13).Value = -46875029
$ws.Cells.Item(139, 14).Value = -22205.6

# CUL row 140
$ws.Cells.Item(140, 8).Value = 1789.625
$ws.Cells.Item(140, 9).Value = 1497.5264
$ws.Cells.Item(140, 11).Value = 4492.5792
$ws.Cells.Item(140, 13).Value = 687.4207999999999

$ws = $wb.Worksheets.Item("GSM")
# GSM row 47
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).ClearContents()

# GSM row 52
$ws.Cells.Item(52, 8).Value = 8666.833000000001
$ws.Cells.Item(52, 10).Value = 8666.833000000001
$ws.Cells.Item(52, 12).Value = 8666.833000000001
$ws.Cells.Item(52, 14).Value = -9184.833000000001

# GSM row 80
$ws.Cells.Item(80, 8).Value = 41685564
$ws.Cells.Item(80, 10).Value = 100022000
$ws.Cells.Item(80, 12).Value = 100022000
$ws.Cells.Item(80, 14).Value = -100023996

# GSM row 83
$ws.Cells.Item(83, 8).Value = 41685564
$ws.Cells.Item(83, 10).Value = 100022000
$ws.Cells.Item(83, 12).Value = 500110000
$ws.Cells.Item(83, 14).Value = -500119984

# GSM row 126
$ws.Cells.Item(126, 8).Value = 2640
$ws.Cells.Item(126, 9).Value = 2698.4
$ws.Cells.Item(126, 10).Value = 2567
$ws.Cells.Item(126, 11).Value = 8095.200000000001
$ws.Cells.Item(126, 12).Value = 7701
$ws.Cells.Item(126, 13).Value = -5625.200000000001
$ws.Cells.Item(126, 14).Value = -12641

# GSM row 134
$ws.Cells.Item(134, 8).Value = 42500
$ws.Cells.Item(134, 10).Value = 42500
$ws.Cells.Item(134, 12).Value = 127500
$ws.Cells.Item(134, 14).Value = -132570

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Cells.Item(22, 8).Value = 4650.769
$ws.Cells.Item(22, 9).Value = 1704.125
$ws.Cells.Item(22, 11).Value = 1704.125
$ws.Cells.Item(22, 13).Value = -1409.125

# LTW row 27
$ws.Cells.Item(27, 8).Value = 4650.769
$ws.Cells.Item(27, 9).Value = 1704.125
$ws.Cells.Item(27, 11).Value = 1704.125
$ws.Cells.Item(27, 13).Value = -1597.125

# LTW row 132
$ws.Cells.Item(132, 8).Value = 3111.1853
$ws.Cells.Item(132, 9).Value = 3375.1765
$ws.Cells.Item(132, 10).Value = 2662.4
$ws.Cells.Item(132, 11).Value = 10125.5295
$ws.Cells.Item(132, 12).Value = 7987.200000000001
$ws.Cells.Item(132, 13).Value = -7595.529500000001
$ws.Cells.Item(132, 14).Value = -13047.2

$ws = $wb.Worksheets.Item("WVR")
# WVR row 54
$ws.Cells.Item(54, 8).Value = 25999.4
$ws.Cells.Item(54, 10).Value = 49998.5
$ws.Cells.Item(54, 12).Value = 49998.5
$ws.Cells.Item(54, 14).Value = -51038.5

# WVR row 136
$ws.Cells.Item(136, 8).Value = 1087.6364
$ws.Cells.Item(136, 9).Value = 849.3570999999999
$ws.Cells.Item(136, 11).Value = 2548.0713
$ws.Cells.Item(136, 13).Value = 1.92870000000039

Write-Output "Applied profit updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"